$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Obstacles to avoid"
$ws.Range("F4").Value = "Changing the position of the target"

$ws.Range("E6").Value = "Obstacles to avoid"
$ws.Range("F6").Value = "Gather data on how it is performing"

$ws.Range("G4").Value = "Write the presentation"

$ws.Columns.Item(5).ColumnWidth = 16.498697916666668
$ws.Columns.Item(6).ColumnWidth = 30.053385416666668
$ws.Columns.Item(7).ColumnWidth = 20.053385416666668

$ws.Range("D4").Select()
